# "Reporte planificación" worksheet:
#  - fix swapped values in G4/J4 and G5/J5 (rows for contract #2 had the
#    "fabricación" duration numbers reversed with the values belonging to
#    the next row)
#  - append a new data row (row 8) for a 5th contract, including a new
#    "12" text label in column M

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte planificación")

# Fix the swapped numbers in rows 4 and 5
$ws.Range("G4").Value = 6
$ws.Range("J4").Value = 7

$ws.Range("G5").Value = 10
$ws.Range("J5").Value = 11

# Add new row 8 with data
$ws.Range("A8").Value = 5

$ws.Range("B8").Value = 42824
$ws.Range("C8").Value = 42830
$ws.Range("B8:C8").NumberFormat = "yyyy-mm-dd"

$ws.Range("D8").Value = 5

$ws.Range("E8").Value = 42831
$ws.Range("F8").Value = 42837
$ws.Range("E8:F8").NumberFormat = "yyyy-mm-dd"

$ws.Range("G8").Value = 2

$ws.Range("H8").Value = 42838
$ws.Range("I8").Value = 42846
$ws.Range("H8:I8").NumberFormat = "yyyy-mm-dd"

$ws.Range("J8").Value = 7

$ws.Range("K8").Value = 42849
$ws.Range("L8").Value = 42852
$ws.Range("K8:L8").NumberFormat = "yyyy-mm-dd"

# M8 holds a new text label ("12"); force it to be stored as text (not a
# number) the same way the other labels in column M are, then restore the
# default (unstyled) cell style so it matches its neighbours.
$m8 = $ws.Range("M8")
$m8.NumberFormat = "@"
$m8.Value = "12"
$m8.Style = "Normal"
